# Automatische test-sync: 2025-06-22 19:06:50
#
# 1. Append a new row (36) to the "Logs" sheet with a new e-mail entry.
# 2. Extend the two conditional-formatting ranges on "Logs" (columns D and G)
#    so they keep covering the data through the new row.
# 3. Update the "Dashboard" sheet's category/count table (rows 9-15) to
#    reflect the re-ranked counts after the new "Juridisch / Contract" entry
#    was added (its count goes from 1 to 2, re-ordering the ties).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Logs sheet - add the new row
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(36, 1).Value = "Privacybeleid"
$logs.Cells.Item(36, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(36, 3).Value = "Wat doen jullie met klantgegevens volgens GDPR?"
$logs.Cells.Item(36, 4).Value = "Juridisch / Contract"
$logs.Cells.Item(36, 5).Value = "Geachte klant,`nDank u wel voor uw vraag over ons beleid met betrekking tot klantgegevens volgens de GDPR. Bij ons bedrijf nemen wij de bescherming van persoonlijke gegevens zeer serieus. Wij voldoen aan de eisen en verplichtingen van de Algemene verordening gegevensbescherming (AVG) om ervoor te zorgen dat alle klantgegevens veilig en vertrouwelijk worden behandeld.`nAls u meer specifieke informatie wenst over hoe wij omgaan met klantgegevens of als u wilt weten welke gegevens wij precies verzamelen en hoe we die gebruiken, dan kunt u contact met ons opnemen. `nMet vriendelijke groet,`n[Bedrijfsnaam] Beveiliging & Compliance Team"
$logs.Cells.Item(36, 6).Value = "2025-06-22 19:06:35"
$logs.Cells.Item(36, 7).Value = "Ja"

# Re-fit the new row's height back to the sheet default (no explicit/custom
# row height), matching how the other data rows are stored.
$logs.Rows.Item(36).AutoFit()

# ---------------------------------------------------------------------
# 2) Logs sheet - widen the conditional formatting ranges to row 36
# ---------------------------------------------------------------------
$dFormats = $logs.Range("D2:D35").FormatConditions
$dFormats.Item(1).ModifyAppliesToRange($logs.Range("D2:D36"))

$gFormats = $logs.Range("G2:G35").FormatConditions
$gFormats.Item(1).ModifyAppliesToRange($logs.Range("G2:G36"))

# ---------------------------------------------------------------------
# 3) Dashboard sheet - refresh the category counts / ordering
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Cells.Item(9, 1).Value = "Juridisch / Contract"
$dash.Cells.Item(9, 2).Value = 2

$dash.Cells.Item(11, 1).Value = "Openingstijden / Locatie"
$dash.Cells.Item(11, 2).Value = 2

$dash.Cells.Item(13, 1).Value = "Klacht / Probleem"
$dash.Cells.Item(13, 2).Value = 1

$dash.Cells.Item(14, 1).Value = "Factuur / Administratie"
$dash.Cells.Item(14, 2).Value = 1

$dash.Cells.Item(15, 1).Value = "Bestelling / Levering"
$dash.Cells.Item(15, 2).Value = 1
